# Apply the cryptos list update (GitHub Actions automated refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing text formatting of the Price column (D) so that
# numeric-looking strings (e.g. "1.00", "37.50") are not coerced into numbers
# and keep their original trailing zeros / grouping dots as plain text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.262.52"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "2.584.08"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "561.54"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "142.78"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "2.591.67"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  +11.39%  "
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Value = "3.038.68"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("E15").Value = "  +6.50%  "
$ws.Range("D16").Value = "59.218.68"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D18").Value = "2.588.34"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "4.59"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "337.43"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "10.38"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "64.15"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").Value = "0.468"
$ws.Range("E25").Value = "  +5.36%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("D28").Value = "7.35"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "0.0₃0776"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").Value = "6.12"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "159.06"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").Value = "19.01"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").Value = "4.05"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "0.873"
$ws.Range("E37").Value = "  -4.26%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "0.877"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "37.50"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "293.55"
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("D43").Value = "132.90"
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("D45").Value = "0.0974"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "0.596"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "10.65"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "19.05"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -0.55%  "

Write-Host "Applied cryptos list update"
